# Daily attendance processing - reorder "Recorded By" entries so that
# real user identifiers (e.g. email addresses) are listed before the
# generic "System"/"system" entries within the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    if ($text -notmatch "backup@backdoor\.com" -and $text -notmatch "dnasr281@gmail\.com") {
        continue
    }

    $parts = $text.Split(",")
    $others = @()
    $systemEntries = @()

    foreach ($p in $parts) {
        $trimmed = $p.Trim()
        if ($trimmed -eq "System" -or $trimmed -eq "system") {
            $systemEntries += $trimmed
        } else {
            $others += $trimmed
        }
    }

    $reordered = @()
    $reordered += $others
    $reordered += $systemEntries

    $newText = [string]::Join(", ", $reordered)

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
